# "procedure code syntax correction"
# Fix the duplicated-word typo "za_procedure_procedure" -> "za_procedure" in the
# Coding System column (C) of the "sterilisation" sheet, widen that column so the
# longer-looking values are readable, and restore the view/selection state that was
# left behind in each sheet (active sheet, selected cell, scroll position).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Data fix: every cell that was tagged with the mistaken "za_procedure_procedure"
#    code should just say "za_procedure" (matching the rest of the column).
# ---------------------------------------------------------------------------
$wsSteril = $wb.Worksheets.Item("sterilisation")

$fixCells = @("C11", "C13", "C15", "C19", "C21", "C52", "C55", "C56", "C59", "C60")
foreach ($cellRef in $fixCells) {
    $wsSteril.Range($cellRef).Value = "za_procedure"
}

# Column C now holds the longer "za_procedure"/"cvv_procedure"/"cbv_procedure"
# strings comfortably - widen it to fit.
$wsSteril.Columns.Item(3).ColumnWidth = 45

# ---------------------------------------------------------------------------
# 2) View/selection bookkeeping restored per sheet.
# ---------------------------------------------------------------------------

# sterilisation: no longer the tab shown when the workbook opens; selection
# moves to the cell that was just corrected.
$wsSteril.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$wsSteril.Range("C11").Select()

# iud: scrolled further down, selection unchanged.
$wsIud = $wb.Worksheets.Item("iud")
$wsIud.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$wsIud.Range("A49").Select()

# preg_test: becomes the active tab, scrolled further down, selection unchanged.
$wsPreg = $wb.Worksheets.Item("preg_test")
$wsPreg.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$wsPreg.Range("C36").Select()
